# Update the "Förändrad" date (column C) for rows 2-27 from 2023-10-22 (45221)
# to 2023-10-25 (45224) on the active worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 27; $row++) {
    $ws.Cells.Item($row, 3).Value = 45224
}
